$wb = $excel.ActiveWorkbook

# ---- Sheet: PIR ----
$ws = $wb.Worksheets.Item("PIR")
$data = @(
    @("2026-02-01","18:37:29","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:29","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:33","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:38","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:43","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:48","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:53","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:37:58","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:38:03","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:38:09","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:38:13","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:38:19","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:38:24","18:00","Bathroom","No Motion","Inactive"),
    @("2026-02-01","18:38:29","18:00","Bathroom","No Motion","Inactive")
)
$startRow = 107
$ws.Range("A107:F120").NumberFormat = "@"
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $data[$i][$col - 1]
    }
}

# ---- Sheet: Humidity ----
$ws = $wb.Worksheets.Item("Humidity")
$data = @(
    @("2026-02-01","18:37:30","18:00","Bathroom","77.9%","Active"),
    @("2026-02-01","18:37:35","18:00","Bathroom","78.9%","Active"),
    @("2026-02-01","18:37:40","18:00","Bathroom","78.0%","Active"),
    @("2026-02-01","18:37:45","18:00","Bathroom","79.0%","Active"),
    @("2026-02-01","18:37:50","18:00","Bathroom","78.2%","Active"),
    @("2026-02-01","18:37:55","18:00","Bathroom","79.2%","Active"),
    @("2026-02-01","18:38:00","18:00","Bathroom","78.2%","Active"),
    @("2026-02-01","18:38:05","18:00","Bathroom","79.2%","Active"),
    @("2026-02-01","18:38:10","18:00","Bathroom","78.3%","Active"),
    @("2026-02-01","18:38:15","18:00","Bathroom","79.3%","Active"),
    @("2026-02-01","18:38:20","18:00","Bathroom","78.3%","Active"),
    @("2026-02-01","18:38:25","18:00","Bathroom","79.3%","Active")
)
$startRow = 160
$ws.Range("A160:F171").NumberFormat = "@"
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $data[$i][$col - 1]
    }
}

# ---- Sheet: Temperature ----
$ws = $wb.Worksheets.Item("Temperature")
$data = @(
    @("2026-02-01","18:37:30","18:00","Bathroom","30.1C","Active"),
    @("2026-02-01","18:37:35","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:37:40","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:37:45","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:37:50","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:37:55","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:38:01","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:38:06","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:38:11","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:38:16","18:00","Bathroom","30.0C","Active"),
    @("2026-02-01","18:38:21","18:00","Bathroom","29.9C","Active"),
    @("2026-02-01","18:38:26","18:00","Bathroom","30.0C","Active")
)
$startRow = 160
$ws.Range("A160:F171").NumberFormat = "@"
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($r, $col).Value = $data[$i][$col - 1]
    }
}

Write-Host "Done updating PIR, Humidity, Temperature sheets"